$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve column D as text so numeric-looking strings keep their exact formatting
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Range("D2").Value = "69.542.20"
$ws.Range("E2").Value = "  +1.97%  "
$ws.Range("D3").Value = "3.381.43"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "581.18"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").Value = "178.80"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("E9").Value = "  +8.18%  "
$ws.Range("D10").Value = "0.588"
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("D11").Value = "48.54"
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("E12").Value = "  +4.06%  "
$ws.Range("D13").Value = "687.97"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").Value = "8.62"
$ws.Range("E14").Value = "  +2.35%  "
$ws.Range("D15").Value = "3.925.00"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").Value = "69.593.82"
$ws.Range("E16").Value = "  +2.04%  "
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").Value = "3.379.54"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("D19").Value = "17.76"
$ws.Range("E19").Value = "  +1.85%  "
$ws.Range("D20").Value = "11.26"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("D22").Value = "17.27"
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("D23").Value = "5.34"
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("D24").Value = "101.49"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("D25").Value = "3.88"
$ws.Range("D26").Value = "2.69"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").Value = "9.71"
$ws.Range("E27").Value = "  +1.89%  "
$ws.Range("D28").Value = "33.60"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("D29").Value = "8.72"
$ws.Range("E29").Value = "  +2.56%  "
$ws.Range("D30").Value = "6.91"
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("D31").Value = "3.86"
$ws.Range("E31").Value = "  +17.58%  "
$ws.Range("D32").Value = "11.06"
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").Value = "554.64"
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("D35").Value = "57.96"
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").Value = "3.604.64"
$ws.Range("E37").Value = "  -2.27%  "
$ws.Range("E38").Value = "  +2.57%  "
$ws.Range("D39").Value = "35.33"
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("D40").Value = "0.0₃0728"
$ws.Range("E40").Value = "  +8.63%  "
$ws.Range("E41").Value = "  +4.17%  "
$ws.Range("D42").Value = "2.73"
$ws.Range("E42").Value = "  +4.39%  "

# Row 43 ("ApeXProtocol") removed; rows 44-51 shift up, THORChain appended as new row 51
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0426"
$ws.Range("E43").Value = "  +2.90%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "0.336"
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "2.66"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "0.129"
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "1.38"
$ws.Range("E48").Value = "  +3.68%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "129.43"
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("B50").Value = "CoreDAO"
$ws.Range("C50").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D50").Value = "2.58"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "7.44"
$ws.Range("E51").Value = "  -0.11%  "

# Restore original (default) style on column D now that values are set
$colD.Style = "Normal"
